$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 46983.953
$ws.Range("I132").Value = 8409.091
$ws.Range("J132").Value = 89416.3
$ws.Range("K132").Value = 25227.273
$ws.Range("L132").Value = 268248.9
$ws.Range("M132").Value = -22697.273
$ws.Range("N132").Value = -273308.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2167.353
$ws.Range("I2").Value = 2205.9375
$ws.Range("J2").Value = 1550
$ws.Range("K2").Value = 2205.9375
$ws.Range("L2").Value = 1550
$ws.Range("M2").Value = -2092.9375
$ws.Range("N2").Value = -1776
$ws.Range("H4").Value = 320
$ws.Range("I4").Value = 280
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 280
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -164
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 255.91667
$ws.Range("I5").Value = 276.8889
$ws.Range("J5").Value = 193
$ws.Range("K5").Value = 276.8889
$ws.Range("L5").Value = 193
$ws.Range("M5").Value = -164.8889
$ws.Range("N5").Value = -417
$ws.Range("H6").Value = 3300001
$ws.Range("I6").Value = 6537502
$ws.Range("J6").Value = 62500
$ws.Range("K6").Value = 6537502
$ws.Range("L6").Value = 62500
$ws.Range("M6").Value = -6537329
$ws.Range("N6").Value = -62846
$ws.Range("H63").Value = 3034.7083
$ws.Range("I63").Value = 2280.5386
$ws.Range("J63").Value = 3926
$ws.Range("K63").Value = 2280.5386
$ws.Range("L63").Value = 3926
$ws.Range("M63").Value = -1594.5386
$ws.Range("N63").Value = -5298
$ws.Range("H66").Value = 3034.7083
$ws.Range("I66").Value = 2280.5386
$ws.Range("J66").Value = 3926
$ws.Range("K66").Value = 11402.693
$ws.Range("L66").Value = 19630
$ws.Range("M66").Value = -7970.692999999999
$ws.Range("N66").Value = -26494
$ws.Range("H88").Value = 12907247
$ws.Range("I88").Value = 25003362
$ws.Range("J88").Value = 3230355.5
$ws.Range("K88").Value = 25003362
$ws.Range("L88").Value = 3230355.5
$ws.Range("M88").Value = -25002956
$ws.Range("N88").Value = -3231167.5
$ws.Range("H91").Value = 12907247
$ws.Range("I91").Value = 25003362
$ws.Range("J91").Value = 3230355.5
$ws.Range("K91").Value = 25003362
$ws.Range("L91").Value = 3230355.5
$ws.Range("M91").Value = -25001958
$ws.Range("N91").Value = -3233163.5
$ws.Range("H116").Value = 2167.353
$ws.Range("I116").Value = 2205.9375
$ws.Range("J116").Value = 1550
$ws.Range("K116").Value = 2205.9375
$ws.Range("L116").Value = 1550
$ws.Range("M116").Value = 88.0625
$ws.Range("N116").Value = -6138
$ws.Range("H132").Value = 16668987
$ws.Range("I132").Value = 23810934
$ws.Range("J132").Value = 4443.6665
$ws.Range("K132").Value = 71432802
$ws.Range("L132").Value = 13330.9995
$ws.Range("M132").Value = -71430272
$ws.Range("N132").Value = -18390.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2167.353
$ws.Range("I3").Value = 2205.9375
$ws.Range("J3").Value = 1550
$ws.Range("K3").Value = 2205.9375
$ws.Range("L3").Value = 1550
$ws.Range("M3").Value = -2091.9375
$ws.Range("N3").Value = -1778
$ws.Range("H4").Value = 255.91667
$ws.Range("I4").Value = 276.8889
$ws.Range("J4").Value = 193
$ws.Range("K4").Value = 276.8889
$ws.Range("L4").Value = 193
$ws.Range("M4").Value = -161.8889
$ws.Range("N4").Value = -423
$ws.Range("H15").Value = 29374.75
$ws.Range("H19").Value = 32330
$ws.Range("J19").Value = 32330
$ws.Range("L19").Value = 32330
$ws.Range("N19").Value = -32676
$ws.Range("H35").Value = 35431
$ws.Range("J35").Value = 35431
$ws.Range("L35").Value = 35431
$ws.Range("N35").Value = -36051
$ws.Range("H82").Value = 12128.5
$ws.Range("I82").Value = 12128.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 12128.5
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -11745.5
$ws.Range("H85").Value = 12128.5
$ws.Range("I85").Value = 12128.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 12128.5
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -10802.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 309.42856
$ws.Range("I7").Value = 280.57144
$ws.Range("J7").Value = 338.2857
$ws.Range("K7").Value = 280.57144
$ws.Range("L7").Value = 338.2857
$ws.Range("M7").Value = -167.57144
$ws.Range("N7").Value = -564.2857
$ws.Range("H31").Value = 5306.3276
$ws.Range("I31").Value = 2627.4
$ws.Range("J31").Value = 5831.608
$ws.Range("K31").Value = 2627.4
$ws.Range("L31").Value = 5831.608
$ws.Range("M31").Value = -2332.4
$ws.Range("N31").Value = -6421.608
$ws.Range("H34").Value = 5306.3276
$ws.Range("I34").Value = 2627.4
$ws.Range("J34").Value = 5831.608
$ws.Range("K34").Value = 2627.4
$ws.Range("L34").Value = 5831.608
$ws.Range("M34").Value = -2425.4
$ws.Range("N34").Value = -6235.608
$ws.Range("H58").Value = 1872.7297
$ws.Range("I58").Value = 1399.8948
$ws.Range("J58").Value = 2371.8333
$ws.Range("K58").Value = 1399.8948
$ws.Range("L58").Value = 2371.8333
$ws.Range("M58").Value = -1196.8948
$ws.Range("N58").Value = -2777.8333
$ws.Range("H107").Value = 522.7143
$ws.Range("I107").Value = 395.32352
$ws.Range("J107").Value = 1064.125
$ws.Range("K107").Value = 395.32352
$ws.Range("L107").Value = 1064.125
$ws.Range("M107").Value = 1524.67648
$ws.Range("N107").Value = -4904.125
$ws.Range("H136").Value = 1872.7297
$ws.Range("I136").Value = 1399.8948
$ws.Range("J136").Value = 2371.8333
$ws.Range("K136").Value = 4199.6844
$ws.Range("L136").Value = 7115.499899999999
$ws.Range("M136").Value = -1649.6844
$ws.Range("N136").Value = -12215.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1166.3334
$ws.Range("I47").Value = 1249.5
$ws.Range("K47").Value = 3748.5
$ws.Range("M47").Value = -3317.5
$ws.Range("H113").Value = 5517.4287
$ws.Range("I113").Value = 14857.429
$ws.Range("J113").Value = 847.4286
$ws.Range("K113").Value = 44572.287
$ws.Range("L113").Value = 2542.2858
$ws.Range("M113").Value = -42402.287
$ws.Range("N113").Value = -6882.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 634154.4399999999
$ws.Range("I2").Value = 926814
$ws.Range("J2").Value = 58.666668
$ws.Range("K2").Value = 926814
$ws.Range("L2").Value = 58.666668
$ws.Range("M2").Value = -926701
$ws.Range("N2").Value = -284.666668
$ws.Range("H43").Value = 4603.4
$ws.Range("I43").Value = 4008.5
$ws.Range("K43").Value = 4008.5
$ws.Range("M43").Value = -3857.5
$ws.Range("H46").Value = 23141.75
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 26970.1
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 26970.1
$ws.Range("M46").Value = -3844
$ws.Range("N46").Value = -27282.1
$ws.Range("H57").Value = 21794.615
$ws.Range("I57").Value = 13000
$ws.Range("J57").Value = 22527.5
$ws.Range("K57").Value = 13000
$ws.Range("L57").Value = 22527.5
$ws.Range("M57").Value = -12180
$ws.Range("N57").Value = -24167.5
$ws.Range("H70").Value = 5172.591
$ws.Range("J70").Value = 4600
$ws.Range("L70").Value = 4600
$ws.Range("N70").Value = -5140
$ws.Range("H73").Value = 5172.591
$ws.Range("J73").Value = 4600
$ws.Range("L73").Value = 4600
$ws.Range("N73").Value = -6472
$ws.Range("H80").Value = 8207.143
$ws.Range("I80").Value = 7500
$ws.Range("J80").Value = 8737.5
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 8737.5
$ws.Range("M80").Value = -6502
$ws.Range("N80").Value = -10733.5
$ws.Range("H83").Value = 8207.143
$ws.Range("I83").Value = 7500
$ws.Range("J83").Value = 8737.5
$ws.Range("K83").Value = 37500
$ws.Range("L83").Value = 43687.5
$ws.Range("M83").Value = -32508
$ws.Range("N83").Value = -53671.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2905
$ws.Range("I61").Value = 3082.4285
$ws.Range("J61").Value = 2698
$ws.Range("K61").Value = 3082.4285
$ws.Range("L61").Value = 2698
$ws.Range("M61").Value = -2880.4285
$ws.Range("N61").Value = -3102
$ws.Range("H113").Value = 2905
$ws.Range("I113").Value = 3082.4285
$ws.Range("J113").Value = 2698
$ws.Range("K113").Value = 3082.4285
$ws.Range("L113").Value = 2698
$ws.Range("M113").Value = -912.4285
$ws.Range("N113").Value = -7038

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1181.3334
$ws.Range("J113").Value = 1100
$ws.Range("L113").Value = 3300
$ws.Range("N113").Value = -7640
